$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @(
    @(1.02, 1.071125148936574, 1.080246609679346, 1.075456385561146, 1.086179600876022, 1.029411991895546, 1.076051327029329, 1.082920544831877, 1.07814290706639, 1.088838113156724, 1.028583383724545),
    @(1.02, 1.073223625446383, 1.082254941624109, 1.07732937088837, 1.08814671069342, 1.029443806110048, 1.077802229427904, 1.084744213159908, 1.079830626585506, 1.090621800355178, 1.029199721004257),
    @(1.02, 1.074577272269718, 1.083550442908062, 1.07853738396962, 1.08941511382647, 1.029461631370708, 1.078930786645533, 1.085919820339161, 1.080918330780363, 1.091771101230529, 1.029596215435513),
    @(1.02, 1.075145363339365, 1.084094131478095, 1.079044312416783, 1.089947307926418, 1.029468464779014, 1.079404201339657, 1.086413007115262, 1.081374577202359, 1.092253123806335, 1.029762354393039),
    @(1.02, 1.075240691473403, 1.084185364718426, 1.079129374590819, 1.090036605167986, 1.029469573432489, 1.079483629916049, 1.086495755073676, 1.081451123493964, 1.092333991067364, 1.029790217974796),
    @(1.02, 1.074584866955018, 1.083557711364275, 1.078544161159237, 1.08942222909637, 1.02946172527281, 1.078937116457479, 1.085926414376253, 1.080924431165407, 1.091777546504628, 1.029598437532878),
    @(1.02, 1.071835225466577, 1.08092618201609, 1.076090197215815, 1.086845329488692, 1.029423315629621, 1.076643975217175, 1.083537792061353, 1.078714195397702, 1.089441938404183, 1.028792162022911),
    @(1.02, 1.066956635245977, 1.076257227713995, 1.071734900698971, 1.082269420302912, 1.02933446178786, 1.072568562223681, 1.079293841371625, 1.074785149936307, 1.085288098339452, 1.027353321257044),
    @(1.02, 1.063680195087045, 1.073121695310586, 1.06880909274052, 1.079193799232727, 1.029260949351283, 1.069827009215032, 1.076439702578613, 1.072141423155226, 1.082491852399139, 1.026381465467675),
    @(1.02, 1.062255395419102, 1.071758213710595, 1.067536594533941, 1.077855767441031, 1.029225720996292, 1.06863374711582, 1.075197628053412, 1.070990593373066, 1.081274338201117, 1.025957539483135),
    @(1.02, 1.061725216688162, 1.071250859555837, 1.06706306343411, 1.077357793381099, 1.029212124243509, 1.068189566032029, 1.074735305876963, 1.070562185674911, 1.080821063340893, 1.025799598593235),
    @(1.02, 1.061838985083004, 1.071359729772601, 1.067164677201117, 1.07746465482975, 1.029215063950585, 1.068284887815334, 1.074834519448803, 1.07065412344527, 1.08091833963573, 1.02583349913613),
    @(1.02, 1.062211590134228, 1.071716294111333, 1.067497470180555, 1.07781462468955, 1.029224607522348, 1.068597050510011, 1.075159432106779, 1.070955200345827, 1.081236891643716, 1.025944493802302),
    @(1.02, 1.06244103827726, 1.071935865523612, 1.06770239901513, 1.078030123365184, 1.029230419838887, 1.068789257519346, 1.075359493446248, 1.071140578674931, 1.081433023998254, 1.026012817937965),
    @(1.02, 1.063774623598292, 1.073212060973959, 1.068893423918094, 1.079282465578499, 1.029263215684981, 1.069906070256474, 1.076522001776526, 1.072217669740039, 1.082572510980263, 1.026409533811851),
    @(1.02, 1.064609496827003, 1.074011018309676, 1.069639003313248, 1.080066329231607, 1.029282877480457, 1.070604952057324, 1.077249530141088, 1.072891655534134, 1.083285463299968, 1.026657545138473),
    @(1.02, 1.065095879815892, 1.074476481325811, 1.070073348271415, 1.080522940810822, 1.029294018184609, 1.071012006241829, 1.077673287441142, 1.07328419525395, 1.083700669604488, 1.026801906896849),
    @(1.02, 1.065261625655572, 1.074635098862301, 1.070221358073554, 1.080678532145579, 1.029297761325048, 1.071150701616317, 1.077817677158698, 1.073417942715268, 1.083842135312539, 1.026851080079329),
    @(1.02, 1.064519983507619, 1.073925355390633, 1.069559065631972, 1.079982290633783, 1.029280801857512, 1.07053003004994, 1.07717153524775, 1.072819403904806, 1.08320903734457, 1.026630966869336),
    @(1.02, 1.062101893565689, 1.071611319777761, 1.067399495136713, 1.077711594276636, 1.029221811303093, 1.0685051528253, 1.07506378015666, 1.070866566806676, 1.081143114823822, 1.025911821851954),
    @(1.02, 1.060576062870318, 1.070151193037953, 1.06603664845321, 1.076278294156156, 1.029181762418862, 1.067226520093741, 1.073732980718428, 1.069633299360345, 1.079838180857294, 1.025456908709603),
    @(1.02, 1.06138546476326, 1.07092573650968, 1.066759605912372, 1.077038656090899, 1.029203273936077, 1.067904878830529, 1.074438999745319, 1.070287601992226, 1.080530529202209, 1.025698331342044),
    @(1.02, 1.064560432528481, 1.073964064491357, 1.069595187693613, 1.080020265920644, 1.029281740754776, 1.070563885907755, 1.077206779628806, 1.072852053112285, 1.083243572941677, 1.026642977368914),
    @(1.02, 1.068221986716782, 1.077468190650534, 1.072864674216216, 1.083456705669959, 1.029359945731908, 1.073626392282351, 1.080395280022129, 1.075805101331091, 1.086366628775588, 1.02772748688846)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($excelRow, $cols[$c]).Value = $rowVals[$c]
    }
}

Write-Host "Updated vm_pu data for 380 kV case"